# Append the new bitcoin-buy entry recorded on 2025-07-20 to the bottom of
# the sheet (row 33), matching the formatting already used for the other
# "MM/DD/YYYY" text-date rows (e.g. A31 "07/13/2025", A32 "07/16/2025"):
# plain text in column A (no special number format / quote-prefix marker)
# and plain numbers in columns B-D.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 33

# Column A: write the date as literal text. A leading apostrophe forces the
# engine to keep it as text instead of auto-parsing it into a date serial
# number; re-applying the neighboring cell's style afterwards clears the
# "number stored as text" quote-prefix marker that the apostrophe entry
# adds, so the resulting cell format matches A32 (the default style).
$ws.Cells.Item($newRow, 1).Value = "'07/20/2025"
$ws.Cells.Item($newRow, 1).Style = $ws.Cells.Item($newRow - 1, 1).Style

# Columns B-D: plain numeric values.
$ws.Cells.Item($newRow, 2).Value = 0.0004220200000000004
$ws.Cells.Item($newRow, 3).Value = 118477.7972607932
$ws.Cells.Item($newRow, 4).Value = 50
